$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values: 1..20 -> 80..99 (offset of +79)
for ($i = 1; $i -le 20; $i++) {
    $ws.Cells.Item($i, 1).Value = 79 + $i
}

# Update the active selection to E5
$ws.Range("E5").Select()
